# Update results with newer values received from server.
$wb = $excel.ActiveWorkbook

# ---- Sheet "2025" ----
$ws = $wb.Worksheets.Item("2025")
$ws.Range("B2").Value = 0.00004227743582363441
$ws.Range("E2").Value = 0.371558911252095
$ws.Range("I2").Value = 0.65827295846499
$ws.Range("L2").Value = 0.3051931082016766
$ws.Range("M2").Value = 0.08616133333333333
$ws.Range("N2").Value = 12.85838940444553
$ws.Range("O2").Value = 3.053012352016873

# ---- Sheet "2030" ----
$ws = $wb.Worksheets.Item("2030")
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = 0.07253249735152435
$ws.Range("E2").Value = 0.3625865573148273
$ws.Range("I2").Value = 0.857380294591765
$ws.Range("M2").Value = 0.08230458333333336
$ws.Range("N2").Value = 9.168045125790171
$ws.Range("O2").Value = 3.528007792201617

# ---- Sheet "2035" ----
$ws = $wb.Worksheets.Item("2035")
$ws.Range("A2").Value = 0.2440261591595064
$ws.Range("B2").Value = 0.0422413368519816
$ws.Range("E2").Value = 0.0572515092927538
$ws.Range("I2").Value = 0.4129214781397788
$ws.Range("M2").Value = 0.0475770833333333
$ws.Range("N2").Value = 3.941306640199873
$ws.Range("O2").Value = 6.977422630444634
